# Add a new title/image pair row to the price list.
#
# The sheet lists game titles in column A paired with an image URL in
# column B. A new row is inserted right before the existing
# "Battlefield 6 Standard Edition PS5" row (row 522), introducing a
# second listing for the same game ("... PS5 НАВСЕГДА") that reuses the
# same product image as the row immediately below it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new blank row at 522 - everything from 522 down
# (through 533) shifts down to 523..534, matching the diff's
# renumbered A523..A534 / B523..B534 cells.
$ws.Rows("522:522").Insert()

# The freshly inserted row inherits the formatting of the row above it
# (the hyperlink-styled B521). The source row (now at 523) carries no
# explicit cell style, so strip formatting to match.
$ws.Range("A522:B522").ClearFormats()

# New title text (goes into shared strings as a brand-new unique entry).
$ws.Range("A522").Value = "Battlefield 6 Standard Edition PS5 НАВСЕГДА"

# Re-use the same product image URL as the row below (the original
# "Battlefield 6 Standard Edition PS5" row, now shifted to 523).
$ws.Range("B522").Value = $ws.Range("B523").Text

# Reflect the author's last on-screen selection when they saved.
[void]$ws.Range("B520").Select()
